$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.258.45'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = '3.304.92'
$ws.Range("E3").Value = '  -2.37%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '189.87'
$ws.Range("E5").Value = '  +2.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '560.33'
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.588'
$ws.Range("E8").Value = '  -2.17%  '
$ws.Range("D9").Value = '3.298.26'
$ws.Range("E9").Value = '  -2.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.185'
$ws.Range("E10").Value = '  -2.39%  '
$ws.Range("E11").Value = '  -1.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.79'
$ws.Range("E12").Value = '  -1.37%  '
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("E14").Value = '  -1.35%  '
$ws.Range("D15").Value = '3.840.39'
$ws.Range("E15").Value = '  -2.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '612.98'
$ws.Range("E16").Value = '  +0.77%  '
$ws.Range("D17").Value = '66.312.74'
$ws.Range("E17").Value = '  -0.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.10'
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.118'
$ws.Range("E19").Value = '  -0.37%  '
$ws.Range("D20").Value = '3.314.06'
$ws.Range("E20").Value = '  -2.08%  '
$ws.Range("E21").Value = '  -5.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.912'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.30'
$ws.Range("E23").Value = '  +7.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.26'
$ws.Range("E24").Value = '  +3.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.09'
$ws.Range("E25").Value = '  -0.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.00'
$ws.Range("E26").Value = '  -2.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.99'
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.75'
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.77'
$ws.Range("E29").Value = '  +2.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.65'
$ws.Range("E30").Value = '  -2.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.28'
$ws.Range("E31").Value = '  -2.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.09'
$ws.Range("E32").Value = '  +3.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.69'
$ws.Range("E33").Value = '  +4.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '565.75'
$ws.Range("E34").Value = '  +2.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.11'
$ws.Range("E35").Value = '  -1.31%  '
$ws.Range("E36").Value = '  -0.67%  '
$ws.Range("D37").Value = '3.752.00'
$ws.Range("E37").Value = '  -3.44%  '
$ws.Range("E38").Value = '  -1.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("D40").Value = '0.0₃0731'
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.32'
$ws.Range("E41").Value = '  -4.44%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '34.04'
$ws.Range("E42").Value = '  +4.66%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.130'
$ws.Range("E43").Value = '  +0.62%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.73'
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.43'
$ws.Range("E45").Value = '  -2.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0427'
$ws.Range("E47").Value = '  +1.08%  '
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("E49").Value = '  -1.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.59'
$ws.Range("E50").Value = '  -5.05%  '
